$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.394.78'
$ws.Range('E2').Value = '  -0.68%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.496.05'
$ws.Range('E3').Value = '  -1.75%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '617.47'
$ws.Range('E5').Value = '  +3.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.86'
$ws.Range('E6').Value = '  -1.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.492.84'
$ws.Range('E7').Value = '  -1.60%  '
$ws.Range('E8').Value = '  -1.67%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('E10').Value = '  +0.88%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.19'
$ws.Range('E11').Value = '  -2.80%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.575'
$ws.Range('E12').Value = '  -1.85%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '45.34'
$ws.Range('E13').Value = '  -2.01%  '
$ws.Range('E14').Value = '  -1.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.058.62'
$ws.Range('E15').Value = '  -1.82%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.32'
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '599.21'
$ws.Range('E17').Value = '  -2.01%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '70.451.12'
$ws.Range('E18').Value = '  -0.55%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.485.19'
$ws.Range('E19').Value = '  -2.26%  '
$ws.Range('E20').Value = '  +1.68%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.47'
$ws.Range('E21').Value = '  +0.82%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.871'
$ws.Range('E22').Value = '  -0.77%  '
$ws.Range('E23').Value = '  -2.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '15.42'
$ws.Range('E24').Value = '  -1.57%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '96.69'
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.67'
$ws.Range('E26').Value = '  -1.00%  '
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('E28').Value = '  -3.25%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.28'
$ws.Range('E29').Value = '  -2.07%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.86'
$ws.Range('E30').Value = '  -1.89%  '
$ws.Range('E31').Value = '  -2.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.90'
$ws.Range('E32').Value = '  -5.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.27'
$ws.Range('E33').Value = '  -1.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.68'
$ws.Range('E34').Value = '  -5.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '593.25'
$ws.Range('E35').Value = '  -13.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0495'
$ws.Range('E36').Value = '  +3.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '10.79'
$ws.Range('E37').Value = '  +0.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0981'
$ws.Range('E38').Value = '  -2.13%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '56.59'
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('E40').Value = '  +0.19%  '
$ws.Range('E41').Value = '  +0.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.30'
$ws.Range('E42').Value = '  -8.53%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.308.86'
$ws.Range('E43').Value = '  -1.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0₃0713'
$ws.Range('E44').Value = '  +2.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.304'
$ws.Range('E45').Value = '  -4.11%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '31.39'
$ws.Range('E46').Value = '  -3.34%  '
$ws.Range('B47').Value = 'ThetaToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.80'
$ws.Range('E47').Value = '  -4.11%  '
$ws.Range('E48').Value = '  -5.60%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.128'
$ws.Range('E49').Value = '  -1.46%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '133.82'
$ws.Range('E50').Value = '  +0.09%  '
$ws.Range('E51').Value = '  -0.02%  '
